$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row text so "diameter" -> "Diameter" (capital D)
# Order of first-use controls the shared-string table ordering, so set
# these in the same order the target workbook's sharedStrings table uses:
# Major Diameter (Min), Major Diameter (Max), Pitch Diameter (Max), Pitch Diameter (Min)
$ws.Range("C1").Value = "Major Diameter (Min)`n"
$ws.Range("B1").Value = "Major Diameter (Max)`n"
$ws.Range("D1").Value = "Pitch Diameter (Max)`n"
$ws.Range("E1").Value = "Pitch Diameter (Min)`n"

# Row height re-wraps to a shorter height now that the header text is shorter
$ws.Rows.Item(1).RowHeight = 27.6

# Selection changes seen in the diff
$ws.Range("E11").Select()
